# TimeSheet.xlsx update: add "Team member 7: " (Sai Naga Sravani Peraka, MT2012122)
# block, duplicated from the "Team member 3" (Satya Deepthi Bhagi) block, and
# retouch a handful of row-heights / the active selection on Sheet2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Duplicate the "Team member 3" block (rows 52-83) into rows 87-118 so that
#    all of the cell styling (date format, wrap-text, bold headers, etc.) is
#    carried over exactly as Excel would when copy/pasting rows.
# ---------------------------------------------------------------------------
$ws.Range("A52").Copy($ws.Range("A87"))          # "Team member N: " banner
$ws.Range("A53:F53").Copy($ws.Range("A88:F88"))  # column headers
$ws.Range("A54:F83").Copy($ws.Range("A89:F118")) # 30 data rows

# ---------------------------------------------------------------------------
# 2. Fix up the text that must differ from the template block.
# ---------------------------------------------------------------------------
$ws.Cells.Item(87, 1).Value = "Team member 7: "

for ($r = 89; $r -le 118; $r++) {
    $ws.Cells.Item($r, 1).Value = "Sai Naga Sravani Peraka"
    $ws.Cells.Item($r, 2).Value = "MT2012122"
}

# Row 116 (copied from row 81) needs different content: a "Formal Meeting(with
# Sir)" entry logged against 1.5 hrs instead of the SVN-repository entry. The
# template's E81 was a wrapped text cell; E116 must become a plain numeric
# (right-aligned, 2-decimal) cell instead, so re-stamp its style first.
$ws.Cells.Item(116, 3).Value = 41305
$ws.Cells.Item(116, 4).Value = "Formal Meeting(with Sir)"
$ws.Range("F1").Copy($ws.Range("E116"))
$ws.Cells.Item(116, 5).Value = 1.5

# Row 118 (copied from row 83) needs the "SVN Repository Setup" / "Created the
# SVN repositories..." pair instead of the "Coding for MiniProject" pair.
$ws.Cells.Item(118, 4).Value = "SVN Repository Setup"
$ws.Cells.Item(118, 5).Value = "Created the SVN repositories for both Ebay and StudentProfile and done the`ninitial SVN commit by creating the dynamic web projects and basic folder structures"

# ---------------------------------------------------------------------------
# 3. Row-height tweaks on both the old and the new block.
# ---------------------------------------------------------------------------
$ws.Rows.Item(54).RowHeight = 30
$ws.Rows.Item(67).RowHeight = 30
$ws.Rows.Item(81).RowHeight = 45
$ws.Rows.Item(83).RowHeight = 30

$ws.Rows.Item(89).RowHeight = 30
$ws.Rows.Item(102).RowHeight = 30
$ws.Rows.Item(118).RowHeight = 45

# ---------------------------------------------------------------------------
# 4. View state: scroll down to the new block and select its last cell.
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 100
$win.ScrollColumn = 1
$ws.Range("C118").Select()
